$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new record as row 218, pushing the existing rows
# (old 218..276) down to (219..277).
$ws.Rows.Item(218).Insert()

$ws.Cells.Item(218, 1).Value = 7
$ws.Cells.Item(218, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(218, 3).Value = "Ñuble"
$ws.Cells.Item(218, 4).Value = 44943
$ws.Cells.Item(218, 5).Value = 16
$ws.Cells.Item(218, 6).Value = 100112032
$ws.Cells.Item(218, 7).Value = "Zapallo italiano"
$ws.Cells.Item(218, 8).Value = "Sin especificar"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 60
$ws.Cells.Item(218, 11).Value = 9000
$ws.Cells.Item(218, 12).Value = 9000
$ws.Cells.Item(218, 13).Value = 9000
$ws.Cells.Item(218, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(218, 15).Value = "Región del Maule"
$ws.Cells.Item(218, 16).Value = 180
$ws.Cells.Item(218, 17).Value = 50
$ws.Cells.Item(218, 18).Value = "Hortaliza"
